# "Generate Report for Archive"
#
# The handoff status text changed from "Ready for handoff" to
# "In Translation" everywhere it is shown: the per-language status
# columns on the Overview sheet (E & F), and the Status column (C) on
# each per-language report sheet (zh-cn, de-de). Narrow those status
# columns afterwards to match the shorter replacement text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $nrows = $used.Rows.Count
    $ncols = $used.Columns.Count
    for ($r = 1; $r -le $nrows; $r++) {
        for ($c = 1; $c -le $ncols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            # Literal string must be the left-hand operand: some cell text
            # (e.g. "True"/"False") resolves to a real Boolean, and a
            # Boolean-left comparison would coerce the right-hand string.
            if ($oldStatus -eq $cell.Text) {
                $cell.Value = $newStatus
            }
        }
    }
}

# Narrow the status columns to match the shorter replacement text.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
